# Natmi following Dr Hou advice
# Recompute LR-pair (Wnt2b -> Fzd4) stats including "M2" as a valid target
# cluster (previously only ECs/FAPs/sCs were targets), expanding the table
# from a 4x3 sender/target grid (12 rows) to a full 4x4 grid (16 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Wnt2b"
$ws.Cells.Item(2, 3).Value = "Fzd4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.675229333333333
$ws.Cells.Item(2, 8).Value = 5.025688
$ws.Cells.Item(2, 9).Value = 0.2721044738138681
$ws.Cells.Item(2, 10).Value = 0.2721044738138681
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 30.38232766666667
$ws.Cells.Item(2, 14).Value = 91.14698300000001
$ws.Cells.Item(2, 15).Value = 0.4410933377331532
$ws.Cells.Item(2, 16).Value = 0.4410933377331531
$ws.Cells.Item(2, 17).Value = 50.89736652214489
$ws.Cells.Item(2, 18).Value = 458.076298699304
$ws.Cells.Item(2, 19).Value = 0.1200234705666825
$ws.Cells.Item(2, 20).Value = 0.1200234705666824

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Wnt2b"
$ws.Cells.Item(3, 3).Value = "Fzd4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.675229333333333
$ws.Cells.Item(3, 8).Value = 5.025688
$ws.Cells.Item(3, 9).Value = 0.2721044738138681
$ws.Cells.Item(3, 10).Value = 0.2721044738138681
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 22.690535
$ws.Cells.Item(3, 14).Value = 68.07160500000001
$ws.Cells.Item(3, 15).Value = 0.3294232070665772
$ws.Cells.Item(3, 16).Value = 0.3294232070665772
$ws.Cells.Item(3, 17).Value = 38.01184982102667
$ws.Cells.Item(3, 18).Value = 342.10664838924
$ws.Cells.Item(3, 19).Value = 0.08963752842092791
$ws.Cells.Item(3, 20).Value = 0.08963752842092791

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Wnt2b"
$ws.Cells.Item(4, 3).Value = "Fzd4"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.675229333333333
$ws.Cells.Item(4, 8).Value = 5.025688
$ws.Cells.Item(4, 9).Value = 0.2721044738138681
$ws.Cells.Item(4, 10).Value = 0.2721044738138681
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.09645133333333333
$ws.Cells.Item(4, 14).Value = 0.289354
$ws.Cells.Item(4, 15).Value = 0.001400289043537939
$ws.Cells.Item(4, 16).Value = 0.001400289043537939
$ws.Cells.Item(4, 17).Value = 0.1615781028391111
$ws.Cells.Item(4, 18).Value = 1.454202925552
$ws.Cells.Item(4, 19).Value = 0.0003810249133792156
$ws.Cells.Item(4, 20).Value = 0.0003810249133792155

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Wnt2b"
$ws.Cells.Item(5, 3).Value = "Fzd4"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.675229333333333
$ws.Cells.Item(5, 8).Value = 5.025688
$ws.Cells.Item(5, 9).Value = 0.2721044738138681
$ws.Cells.Item(5, 10).Value = 0.2721044738138681
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 15.71027466666667
$ws.Cells.Item(5, 14).Value = 47.130824
$ws.Cells.Item(5, 15).Value = 0.2280831661567317
$ws.Cells.Item(5, 16).Value = 0.2280831661567317
$ws.Cells.Item(5, 17).Value = 26.31831295632355
$ws.Cells.Item(5, 18).Value = 236.864816606912
$ws.Cells.Item(5, 19).Value = 0.06206244991287852
$ws.Cells.Item(5, 20).Value = 0.06206244991287852

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Wnt2b"
$ws.Cells.Item(6, 3).Value = "Fzd4"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.433537333333333
$ws.Cells.Item(6, 8).Value = 7.300612
$ws.Cells.Item(6, 9).Value = 0.3952750721451891
$ws.Cells.Item(6, 10).Value = 0.3952750721451891
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 30.38232766666667
$ws.Cells.Item(6, 14).Value = 91.14698300000001
$ws.Cells.Item(6, 15).Value = 0.4410933377331532
$ws.Cells.Item(6, 16).Value = 0.4410933377331531
$ws.Cells.Item(6, 17).Value = 73.93652865039955
$ws.Cells.Item(6, 18).Value = 665.428757853596
$ws.Cells.Item(6, 19).Value = 0.1743532008952344
$ws.Cells.Item(6, 20).Value = 0.1743532008952344

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Wnt2b"
$ws.Cells.Item(7, 3).Value = "Fzd4"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.433537333333333
$ws.Cells.Item(7, 8).Value = 7.300612
$ws.Cells.Item(7, 9).Value = 0.3952750721451891
$ws.Cells.Item(7, 10).Value = 0.3952750721451891
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 22.690535
$ws.Cells.Item(7, 14).Value = 68.07160500000001
$ws.Cells.Item(7, 15).Value = 0.3294232070665772
$ws.Cells.Item(7, 16).Value = 0.3294232070665772
$ws.Cells.Item(7, 17).Value = 55.21826403580666
$ws.Cells.Item(7, 18).Value = 496.96437632226
$ws.Cells.Item(7, 19).Value = 0.1302127819395409
$ws.Cells.Item(7, 20).Value = 0.1302127819395409

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Wnt2b"
$ws.Cells.Item(8, 3).Value = "Fzd4"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.433537333333333
$ws.Cells.Item(8, 8).Value = 7.300612
$ws.Cells.Item(8, 9).Value = 0.3952750721451891
$ws.Cells.Item(8, 10).Value = 0.3952750721451891
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.09645133333333333
$ws.Cells.Item(8, 14).Value = 0.289354
$ws.Cells.Item(8, 15).Value = 0.001400289043537939
$ws.Cells.Item(8, 16).Value = 0.001400289043537939
$ws.Cells.Item(8, 17).Value = 0.2347179205164444
$ws.Cells.Item(8, 18).Value = 2.112461284648
$ws.Cells.Item(8, 19).Value = 0.0005534993527085768
$ws.Cells.Item(8, 20).Value = 0.0005534993527085767

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Wnt2b"
$ws.Cells.Item(9, 3).Value = "Fzd4"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.433537333333333
$ws.Cells.Item(9, 8).Value = 7.300612
$ws.Cells.Item(9, 9).Value = 0.3952750721451891
$ws.Cells.Item(9, 10).Value = 0.3952750721451891
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 15.71027466666667
$ws.Cells.Item(9, 14).Value = 47.130824
$ws.Cells.Item(9, 15).Value = 0.2280831661567317
$ws.Cells.Item(9, 16).Value = 0.2280831661567317
$ws.Cells.Item(9, 17).Value = 38.23153991825421
$ws.Cells.Item(9, 18).Value = 344.083859264288
$ws.Cells.Item(9, 19).Value = 0.09015558995770528
$ws.Cells.Item(9, 20).Value = 0.09015558995770528

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Wnt2b"
$ws.Cells.Item(10, 3).Value = "Fzd4"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.9268273333333333
$ws.Cells.Item(10, 8).Value = 2.780482
$ws.Cells.Item(10, 9).Value = 0.1505428891644152
$ws.Cells.Item(10, 10).Value = 0.1505428891644152
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 30.38232766666667
$ws.Cells.Item(10, 14).Value = 91.14698300000001
$ws.Cells.Item(10, 15).Value = 0.4410933377331532
$ws.Cells.Item(10, 16).Value = 0.4410933377331531
$ws.Cells.Item(10, 17).Value = 28.15917173175622
$ws.Cells.Item(10, 18).Value = 253.432545585806
$ws.Cells.Item(10, 19).Value = 0.06640346545352406
$ws.Cells.Item(10, 20).Value = 0.06640346545352405

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Wnt2b"
$ws.Cells.Item(11, 3).Value = "Fzd4"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.9268273333333333
$ws.Cells.Item(11, 8).Value = 2.780482
$ws.Cells.Item(11, 9).Value = 0.1505428891644152
$ws.Cells.Item(11, 10).Value = 0.1505428891644152
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 22.690535
$ws.Cells.Item(11, 14).Value = 68.07160500000001
$ws.Cells.Item(11, 15).Value = 0.3294232070665772
$ws.Cells.Item(11, 16).Value = 0.3294232070665772
$ws.Cells.Item(11, 17).Value = 21.03020804595667
$ws.Cells.Item(11, 18).Value = 189.27187241361
$ws.Cells.Item(11, 19).Value = 0.04959232134960995
$ws.Cells.Item(11, 20).Value = 0.04959232134960995

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Wnt2b"
$ws.Cells.Item(12, 3).Value = "Fzd4"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.9268273333333333
$ws.Cells.Item(12, 8).Value = 2.780482
$ws.Cells.Item(12, 9).Value = 0.1505428891644152
$ws.Cells.Item(12, 10).Value = 0.1505428891644152
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.09645133333333333
$ws.Cells.Item(12, 14).Value = 0.289354
$ws.Cells.Item(12, 15).Value = 0.001400289043537939
$ws.Cells.Item(12, 16).Value = 0.001400289043537939
$ws.Cells.Item(12, 17).Value = 0.08939373206977778
$ws.Cells.Item(12, 18).Value = 0.804543588628
$ws.Cells.Item(12, 19).Value = 0.000210803558279477
$ws.Cells.Item(12, 20).Value = 0.000210803558279477

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Wnt2b"
$ws.Cells.Item(13, 3).Value = "Fzd4"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.9268273333333333
$ws.Cells.Item(13, 8).Value = 2.780482
$ws.Cells.Item(13, 9).Value = 0.1505428891644152
$ws.Cells.Item(13, 10).Value = 0.1505428891644152
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 15.71027466666667
$ws.Cells.Item(13, 14).Value = 47.130824
$ws.Cells.Item(13, 15).Value = 0.2280831661567317
$ws.Cells.Item(13, 16).Value = 0.2280831661567317
$ws.Cells.Item(13, 17).Value = 14.56071197524089
$ws.Cells.Item(13, 18).Value = 131.046407777168
$ws.Cells.Item(13, 19).Value = 0.03433629880300176
$ws.Cells.Item(13, 20).Value = 0.03433629880300176

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Wnt2b"
$ws.Cells.Item(14, 3).Value = "Fzd4"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.120972666666667
$ws.Cells.Item(14, 8).Value = 3.362918
$ws.Cells.Item(14, 9).Value = 0.1820775648765275
$ws.Cells.Item(14, 10).Value = 0.1820775648765275
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 30.38232766666667
$ws.Cells.Item(14, 14).Value = 91.14698300000001
$ws.Cells.Item(14, 15).Value = 0.4410933377331532
$ws.Cells.Item(14, 16).Value = 0.4410933377331531
$ws.Cells.Item(14, 17).Value = 34.05775886404378
$ws.Cells.Item(14, 18).Value = 306.519829776394
$ws.Cells.Item(14, 19).Value = 0.08031320081771225
$ws.Cells.Item(14, 20).Value = 0.08031320081771225

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Wnt2b"
$ws.Cells.Item(15, 3).Value = "Fzd4"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1.120972666666667
$ws.Cells.Item(15, 8).Value = 3.362918
$ws.Cells.Item(15, 9).Value = 0.1820775648765275
$ws.Cells.Item(15, 10).Value = 0.1820775648765275
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 22.690535
$ws.Cells.Item(15, 14).Value = 68.07160500000001
$ws.Cells.Item(15, 15).Value = 0.3294232070665772
$ws.Cells.Item(15, 16).Value = 0.3294232070665772
$ws.Cells.Item(15, 17).Value = 25.43546952704333
$ws.Cells.Item(15, 18).Value = 228.91922574339
$ws.Cells.Item(15, 19).Value = 0.05998057535649846
$ws.Cells.Item(15, 20).Value = 0.05998057535649848

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Wnt2b"
$ws.Cells.Item(16, 3).Value = "Fzd4"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.120972666666667
$ws.Cells.Item(16, 8).Value = 3.362918
$ws.Cells.Item(16, 9).Value = 0.1820775648765275
$ws.Cells.Item(16, 10).Value = 0.1820775648765275
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.09645133333333333
$ws.Cells.Item(16, 14).Value = 0.289354
$ws.Cells.Item(16, 15).Value = 0.001400289043537939
$ws.Cells.Item(16, 16).Value = 0.001400289043537939
$ws.Cells.Item(16, 17).Value = 0.1081193083302222
$ws.Cells.Item(16, 18).Value = 0.973073774972
$ws.Cells.Item(16, 19).Value = 0.0002549612191706697
$ws.Cells.Item(16, 20).Value = 0.0002549612191706697

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Wnt2b"
$ws.Cells.Item(17, 3).Value = "Fzd4"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 1.120972666666667
$ws.Cells.Item(17, 8).Value = 3.362918
$ws.Cells.Item(17, 9).Value = 0.1820775648765275
$ws.Cells.Item(17, 10).Value = 0.1820775648765275
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 15.71027466666667
$ws.Cells.Item(17, 14).Value = 47.130824
$ws.Cells.Item(17, 15).Value = 0.2280831661567317
$ws.Cells.Item(17, 16).Value = 0.2280831661567317
$ws.Cells.Item(17, 17).Value = 17.61078848715911
$ws.Cells.Item(17, 18).Value = 158.497096384432
$ws.Cells.Item(17, 19).Value = 0.04152882748314611
$ws.Cells.Item(17, 20).Value = 0.04152882748314612

